# Refatorando artefatos de acordo com o feedback da ac4
#
# Applies the edits made to the "Analise de Eventos" worksheet:
#   - D11: "loja recebe o boleto" -> "Loja recebe o boleto"
#   - D15: "Inicia  transação com a administradora de cartões"
#          -> "Loja inicia transação com a administradora de cartões"
#   - D16: "administradora de cartões devolve transação"
#          -> "Administradora de cartões devolve transação"
#   - the "x(3)" marker moves from G15 to E15
#   - selection/active cell moves to E15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analise de Eventos")
$ws.Activate()

$ws.Range("D11").Value = "Loja recebe o boleto"
$ws.Range("D15").Value = "Loja inicia transação com a administradora de cartões"
$ws.Range("D16").Value = "Administradora de cartões devolve transação"

$ws.Range("G15").Value = ""
$ws.Range("E15").Value = "x(3)"

$ws.Range("E15").Select()
